$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily auto-push: insert today's (2026/02/14) reading as a new row right
# after the existing 2026/02/14 block (row 817), shifting every later row
# down by one. This mirrors the source sheet's append pattern where each
# day's collected samples are appended chronologically before the
# historical forward-looking rows.
$ws.Rows.Item(817).Insert()

# Excel would otherwise auto-detect "2026/02/14" as a date literal and
# reformat/convert the cell; force it to remain literal text like the
# other date column cells, then restore the default (unstyled) look so no
# stray number-format style sticks around on the new row.
$ws.Range("A817").NumberFormat = "@"
$ws.Range("A817").Value = "2026/02/14"
$ws.Range("A817").Style = "Normal"

$ws.Range("B817").Value = "土"
$ws.Range("C817").Value = 7
$ws.Range("D817").Value = 201
